# Commit: "Move loader into a separate application"
# Adds completion dates (27-Dec-2023) for the "testing the screen", "Bugs"
# and "split the utility into a separate application " rows on the Tasks
# sheet, and updates the saved view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$dateSerial = 45287  # 27-Dec-2023 as an Excel date serial

# Use an already-styled date cell (B22, style matches the one needed on
# rows 26/27/29) as the format donor so the existing style is reused
# instead of new cellXfs entries being minted.
$formatDonor = $ws.Range("B22")

# Row 26: "testing the screen "
$ws.Range("B26").Value = $dateSerial
$ws.Range("C26").Value = $dateSerial
$formatDonor.Copy()
$ws.Range("B26:C26").PasteSpecial(-4122) # xlPasteFormats

# Row 27: "Bugs"
$ws.Range("B27").Value = $dateSerial
$ws.Range("C27").Value = $dateSerial
$formatDonor.Copy()
$ws.Range("B27:C27").PasteSpecial(-4122) # xlPasteFormats

# Row 29: "split the utility into a separate application "
$ws.Range("B29").Value = $dateSerial
$ws.Range("C29").Value = $dateSerial
$formatDonor.Copy()
$ws.Range("B29:C29").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# Update the view: scroll position and active selection, matching the
# saved sheetView state after the edits were made.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B30").Select()
